{"js": "const body = context.document.body;\n\n{\n  const results = body.search(\"2025-10-18 Saturday\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"2025-10-19 Sunday\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"59\u00d780=4720\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"87\u00d782=7134\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"93\u00d743=3999\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"71\u00d773=5183\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"96\u00d794=9024\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"90\u00d729=2610\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"93\u00d735=3255\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"12\u00d755=660\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"34\u00d750=1700\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"70\u00d788=6160\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"87\u00d756=4872\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"42\u00d784=3528\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"73\u00d794=6862\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"73\u00d730=2190\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"37\u00d713=481\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"94\u00d727=2538\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"60\u00d712=720\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"27\u00d791=2457\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"83\u00d723=1909\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"26\u00d752=1352\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"16\u00d723=368\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"15\u00d723=345\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"11\u00d796=1056\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"70\u00d798=6860\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"38\u00d750=1900\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"45\u00d739=1755\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"83\u00d721=1743\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"58\u00d798=5684\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"47\u00d772=3384\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"48\u00d785=4080\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"74\u00d721=1554\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"65\u00d732=2080\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"99\u00d776=7524\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"84\u00d739=3276\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"58\u00d742=2436\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"93\u00d720=1860\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"93\u00d775=6975\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"86\u00d744=3784\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"30\u00d740=1200\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"31\u00d755=1705\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"78\u00d793=7254\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"12\u00d778=936\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"90\u00d756=5040\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"71\u00d724=1704\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"81\u00d724=1944\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"53\u00d758=3074\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"58\u00d713=754\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"50\u00d761=3050\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n{\n  const results = body.search(\"88\u00d770=6160\", {matchCase: true});\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"85\u00d760=5100\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"2025-10-18 Saturday\"\n$find.Replacement.Text = \"2025-10-19 Sunday\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"59\u00d780=4720\"\n$find.Replacement.Text = \"87\u00d782=7134\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"93\u00d743=3999\"\n$find.Replacement.Text = \"71\u00d773=5183\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"96\u00d794=9024\"\n$find.Replacement.Text = \"90\u00d729=2610\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"93\u00d735=3255\"\n$find.Replacement.Text = \"12\u00d755=660\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"34\u00d750=1700\"\n$find.Replacement.Text = \"70\u00d788=6160\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"87\u00d756=4872\"\n$find.Replacement.Text = \"42\u00d784=3528\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"73\u00d794=6862\"\n$find.Replacement.Text = \"73\u00d730=2190\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"37\u00d713=481\"\n$find.Replacement.Text = \"94\u00d727=2538\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"60\u00d712=720\"\n$find.Replacement.Text = \"27\u00d791=2457\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"83\u00d723=1909\"\n$find.Replacement.Text = \"26\u00d752=1352\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"16\u00d723=368\"\n$find.Replacement.Text = \"15\u00d723=345\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"11\u00d796=1056\"\n$find.Replacement.Text = \"70\u00d798=6860\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"38\u00d750=1900\"\n$find.Replacement.Text = \"45\u00d739=1755\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"83\u00d721=1743\"\n$find.Replacement.Text = \"58\u00d798=5684\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"47\u00d772=3384\"\n$find.Replacement.Text = \"48\u00d785=4080\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"74\u00d721=1554\"\n$find.Replacement.Text = \"65\u00d732=2080\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"99\u00d776=7524\"\n$find.Replacement.Text = \"84\u00d739=3276\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"58\u00d742=2436\"\n$find.Replacement.Text = \"93\u00d720=1860\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"93\u00d775=6975\"\n$find.Replacement.Text = \"86\u00d744=3784\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"30\u00d740=1200\"\n$find.Replacement.Text = \"31\u00d755=1705\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"78\u00d793=7254\"\n$find.Replacement.Text = \"12\u00d778=936\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"90\u00d756=5040\"\n$find.Replacement.Text = \"71\u00d724=1704\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"81\u00d724=1944\"\n$find.Replacement.Text = \"53\u00d758=3074\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"58\u00d713=754\"\n$find.Replacement.Text = \"50\u00d761=3050\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"88\u00d770=6160\"\n$find.Replacement.Text = \"85\u00d760=5100\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2, $false, $false, $false, $false) | Out-Null\n"}
